$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking text values (e.g. "592.76") that Excel
# would otherwise auto-convert to real numbers. Force the Price column to
# Text format before writing, then restore the default "Normal" style so no
# residual formatting diff is left behind on cells whose value we did not touch.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '63.529.84'
$ws.Range("E2").Value = '  -1.95%  '

$ws.Range("D3").Value = '3.066.54'
$ws.Range("E3").Value = '  -2.67%  '

$ws.Range("E4").Value = '  -0.59%  '

$ws.Range("D5").Value = '592.76'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").Value = '155.19'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").Value = '3.065.24'
$ws.Range("E9").Value = '  -2.67%  '

$ws.Range("E10").Value = '  -3.61%  '

$ws.Range("D11").Value = '5.88'
$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("E12").Value = '  -3.59%  '

$ws.Range("E13").Value = '  -4.89%  '

$ws.Range("D14").Value = '36.52'
$ws.Range("E14").Value = '  -5.35%  '

$ws.Range("E15").Value = '  +0.25%  '

$ws.Range("D16").Value = '3.573.91'
$ws.Range("E16").Value = '  -2.63%  '

$ws.Range("E17").Value = '  -2.14%  '

$ws.Range("D18").Value = '63.403.24'
$ws.Range("E18").Value = '  -1.48%  '

$ws.Range("D19").Value = '3.068.82'
$ws.Range("E19").Value = '  -2.52%  '

$ws.Range("D20").Value = '478.06'
$ws.Range("E20").Value = '  +0.54%  '

$ws.Range("E21").Value = '  -4.37%  '

$ws.Range("E22").Value = '  -5.69%  '

$ws.Range("E23").Value = '  -2.47%  '

$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +0.45%  '

$ws.Range("D25").Value = '81.26'
$ws.Range("E25").Value = '  -1.57%  '

$ws.Range("D26").Value = '12.78'
$ws.Range("E26").Value = '  -5.42%  '

$ws.Range("D27").Value = '10.56'
$ws.Range("E27").Value = '  +6.54%  '

$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.05%  '

$ws.Range("E29").Value = '  +1.34%  '

$ws.Range("E30").Value = '  -1.94%  '

$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("E32").Value = '  -2.83%  '

$ws.Range("E33").Value = '  -5.73%  '

$ws.Range("E34").Value = '  -2.57%  '

$ws.Range("D35").Value = '0.0₃0833'
$ws.Range("E35").Value = '  -5.58%  '

$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("E37").Value = '  -4.19%  '

$ws.Range("E38").Value = '  -5.96%  '

$ws.Range("E39").Value = '  -4.23%  '

$ws.Range("D40").Value = '50.69'
$ws.Range("E40").Value = '  -1.40%  '

$ws.Range("E41").Value = '  -2.33%  '

$ws.Range("D42").Value = '436.72'
$ws.Range("E42").Value = '  -6.71%  '

$ws.Range("D43").Value = '0.289'
$ws.Range("E43").Value = '  -4.44%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = '0.111'
$ws.Range("E44").Value = '  +0.75%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0361'
$ws.Range("E45").Value = '  -5.20%  '

$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").Value = '39.94'
$ws.Range("E46").Value = '  +1.09%  '

$ws.Range("D47").Value = '2.819.69'
$ws.Range("E47").Value = '  -2.78%  '

$ws.Range("D48").Value = '131.87'
$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").Value = '25.63'
$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("E51").Value = '  -3.24%  '

# Restore the column's style to the workbook default now that every literal
# has been committed as text, so the saved file carries no stray "@" format.
$ws.Range("D2:D51").Style = "Normal"
